$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.525.56"
$ws.Range("E2").Value = "  +5.95%  "
$ws.Range("D3").Value = "1.936.98"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'251.31"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'0.691"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'48.53"
$ws.Range("E8").Value = "  +12.23%  "
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("D10").Value = "'58.71"
$ws.Range("E10").Value = "  +7.32%  "
$ws.Range("D11").Value = "'0.0774"
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "'15.74"
$ws.Range("E13").Value = "  +13.08%  "
$ws.Range("E14").Value = "  +8.26%  "
$ws.Range("D15").Value = "2.221.86"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "1.935.13"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "37.673.41"
$ws.Range("E18").Value = "  +6.28%  "
$ws.Range("D19").Value = "'75.44"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("D21").Value = "'13.75"
$ws.Range("E21").Value = "  +6.92%  "
$ws.Range("D22").Value = "'253.44"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").Value = "'5.24"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -5.68%  "
$ws.Range("D26").Value = "'168.88"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'2.15"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'8.96"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").Value = "'18.91"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "'4.60"
$ws.Range("E31").Value = "  +7.14%  "
$ws.Range("D32").Value = "'0.0616"
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").Value = "'0.0917"
$ws.Range("E33").Value = "  +27.30%  "
$ws.Range("D34").Value = "'4.35"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "'19.17"
$ws.Range("E37").Value = "  +38.70%  "
$ws.Range("D38").Value = "'0.899"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").Value = "'105.77"
$ws.Range("E41").Value = "  +7.40%  "
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("D43").Value = "'17.58"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("D44").Value = "'2.89"
$ws.Range("E44").Value = "  +20.51%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "1.353.99"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'0.0844"
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("D49").Value = "'2.81"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'3.84"
$ws.Range("E50").Value = "  +14.92%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'6.45"
$ws.Range("E51").Value = "  +2.30%  "
